$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update country code in row 2 from "PL" to "EL"
$ws.Range("A2").Value = "EL"
